# "Actualizacion Datos Personales 4 nov"
# Adds the missing contact/tutor details for CESAR CUEVAS CUATRA (row 6) on
# sheet "3APM": Correo, Tel_Movil, Tel_Fijo, Tutor, Correo_Tutor y
# Telefono_Tutor.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3APM")

# E6 - Correo
$ws.Range("E6").Value = "cesarcuevasc3@gmail.com"

# F6 - Tel_Movil (keep it text, like every other phone number in the sheet)
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "2722848082"
$ws.Range("F6").Style = "Normal"

# G6 - Tel_Fijo
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "2722848082"
$ws.Range("G6").Style = "Normal"

# H6 - Tutor
$ws.Range("H6").Value = "GERARDO CUEVAS MACUIXTLE"

# I6 - Correo_Tutor
$ws.Range("I6").Value = "cesarcuevasc3@gmail.com"

# J6 - Telefono_Tutor
$ws.Range("J6").NumberFormat = "@"
$ws.Range("J6").Value = "2722848082"
$ws.Range("J6").Style = "Normal"
